$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D): values that look numeric must be forced to text
# using a leading apostrophe, then the style is reset to Normal so the
# cell keeps its original (unstyled) formatting while the stored value
# stays an exact text string identical to the source data feed.
$priceUpdates = @{
    "D2" = "29.349.00"
    "D3" = "1.878.02"
    "D4" = "1.001"
    "D5" = "0.7214"
    "D6" = "242.49"
    "D7" = "1.001"
    "D8" = "0.08022"
    "D9" = "0.3147"
    "D10" = "24.99"
    "D11" = "0.08203"
    "D12" = "1.870.25"
    "D13" = "94.61"
    "D15" = "0.7130"
    "D16" = "6.417"
    "D17" = "0.000008521"
    "D18" = "29.348.03"
    "D19" = "243.45"
    "D20" = "13.30"
    "D21" = "1.001"
    "D22" = "7.750"
    "D23" = "1.002"
    "D24" = "0.1596"
    "D25" = "9.039"
    "D26" = "162.40"
    "D28" = "1.503"
    "D29" = "4.407"
    "D30" = "4.303"
    "D31" = "1.210"
    "D32" = "0.05361"
    "D33" = "1.934"
    "D34" = "0.7658"
    "D35" = "1.177"
    "D36" = "2.710"
    "D37" = "0.01876"
    "D38" = "1.277.74"
    "D39" = "2.750"
    "D40" = "6.437"
    "D41" = "112.89"
    "D42" = "0.9119"
    "D43" = "74.15"
    "D45" = "1.001"
    "D46" = "2.022.24"
    "D47" = "0.5226"
    "D48" = "1.799"
    "D50" = "0.4342"
    "D51" = "7.099"
}
foreach ($cellRef in $priceUpdates.Keys) {
    $ws.Range($cellRef).Value = "'" + $priceUpdates[$cellRef]
    $ws.Range($cellRef).Style = "Normal"
}

# Volume(1h) column (E): plain percentage text, safe to assign directly
$volumeUpdates = @{
    "E2" = "  +0.24%  "
    "E3" = "  +0.35%  "
    "E4" = "  +0.10%  "
    "E5" = "  +1.83%  "
    "E6" = "  +0.39%  "
    "E7" = "  +0.09%  "
    "E8" = "  +2.74%  "
    "E9" = "  +1.79%  "
    "E10" = "  +0.03%  "
    "E11" = "  -2.39%  "
    "E12" = "  -0.34%  "
    "E13" = "  +3.95%  "
    "E14" = "  +0.09%  "
    "E15" = "  +0.41%  "
    "E16" = "  +5.78%  "
    "E17" = "  +4.07%  "
    "E18" = "  +0.20%  "
    "E19" = "  +1.62%  "
    "E20" = "  +0.86%  "
    "E21" = "  +0.08%  "
    "E22" = "  +0.15%  "
    "E23" = "  +0.12%  "
    "E24" = "  +0.43%  "
    "E25" = "  +0.49%  "
    "E26" = "  -0.20%  "
    "E27" = "  +0.44%  "
    "E28" = "  -0.02%  "
    "E29" = "  +0.51%  "
    "E30" = "  +0.21%  "
    "E31" = "  -6.55%  "
    "E32" = "  -0.01%  "
    "E33" = "  -0.53%  "
    "E34" = "  +2.44%  "
    "E35" = "  +0.18%  "
    "E36" = "  +0.49%  "
    "E37" = "  +0.29%  "
    "E38" = "  +4.51%  "
    "E39" = "  +0.96%  "
    "E40" = "  -0.81%  "
    "E41" = "  +4.20%  "
    "E42" = "  +2.38%  "
    "E43" = "  +2.62%  "
    "E44" = "  +6.94%  "
    "E45" = "  +0.11%  "
    "E46" = "  +0.03%  "
    "E47" = "  +0.58%  "
    "E48" = "  +0.35%  "
    "E50" = "  +0.79%  "
    "E51" = "  +0.47%  "
}
foreach ($cellRef in $volumeUpdates.Keys) {
    $ws.Range($cellRef).Value = $volumeUpdates[$cellRef]
}
